$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '255.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.59%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.16'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.322'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '3.13%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05852'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.91%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.704'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.60%'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8664'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.93%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9070'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '5.09%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1420'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.86%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07150'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.05%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03180'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.50%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09221'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.48%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001557'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.91%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'One'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006085'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.31%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005813'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.52%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.499'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.02%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.228'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.43%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.201'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.48%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.57%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03446'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.17%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.19%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.542'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.45%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04157'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.37%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1378'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.17%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.005040'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '21.67%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.001227'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.14%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001200'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.01%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001936'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '33.72%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03851'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.26%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1103'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.95%'
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003833'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-32.61%'
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002343'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.34%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01089'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '14.00%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005222'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.60%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.03%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.08747'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002155'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.24%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.03%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.03%'
